$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 104.017544327471
$ws.Range("C2").Value = 1.327233191755586
$ws.Range("D2").Value = 3.231523654476329
$ws.Range("E2").Value = 9.260013067703733
